# Apply Horarios Linea 141 schedule update (commit: Horarios actualizados Linea 141 - 221)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 15:51:40"
$ws.Range("A3").Value = "Total filas: 355"
$ws.Range("A82").Value = "08:39:56"
$ws.Range("B82").Value = "08:42"
$ws.Range("C82").Value = "14_ABASTO"
$ws.Range("D82").Value = 3
$ws.Range("E82").Value = "LP1912"
$ws.Range("A83").Value = "06:52:23"
$ws.Range("B83").Value = "08:42"
$ws.Range("C83").Value = "81_EL PELIGRO"
$ws.Range("D83").Value = 110
$ws.Range("E83").Value = "LP1912"
$ws.Range("A116").Value = "07:46:15"
$ws.Range("B116").Value = "09:33"
$ws.Range("C116").Value = "10_OLMOS"
$ws.Range("D116").Value = 107
$ws.Range("E116").Value = "LP1912"
$ws.Range("A117").Value = "08:50:00"
$ws.Range("B117").Value = "09:33"
$ws.Range("C117").Value = "16_SANTA ANA"
$ws.Range("D117").Value = 43
$ws.Range("E117").Value = "LP1912"
$ws.Range("A159").Value = "10:28:12"
$ws.Range("B159").Value = "11:04"
$ws.Range("C159").Value = "11_ETCHEVERRY"
$ws.Range("D159").Value = 36
$ws.Range("E159").Value = "LP1912"
$ws.Range("A160").Value = "10:57:58"
$ws.Range("B160").Value = "11:04"
$ws.Range("C160").Value = "23_HERNANDEZ"
$ws.Range("D160").Value = 7
$ws.Range("E160").Value = "LP1912"
$ws.Range("A171").Value = "11:23:54"
$ws.Range("B171").Value = "11:23"
$ws.Range("C171").Value = "17_ROMERO"
$ws.Range("D171").Value = 0
$ws.Range("E171").Value = "LP1912"
$ws.Range("A172").Value = "11:23:54"
$ws.Range("B172").Value = "11:23"
$ws.Range("C172").Value = "16_SANTA ANA"
$ws.Range("D172").Value = 0
$ws.Range("E172").Value = "LP1912"
$ws.Range("A214").Value = "10:28:12"
$ws.Range("B214").Value = "12:21"
$ws.Range("C214").Value = "215A_EL PATO"
$ws.Range("D214").Value = 113
$ws.Range("E214").Value = "LP1912"
$ws.Range("A215").Value = "10:28:12"
$ws.Range("B215").Value = "12:21"
$ws.Range("C215").Value = "26_HERNANDEZ"
$ws.Range("D215").Value = 113
$ws.Range("E215").Value = "LP1912"
$ws.Range("A216").Value = "12:16:51"
$ws.Range("B216").Value = "12:21"
$ws.Range("C216").Value = "16_SANTA ANA"
$ws.Range("D216").Value = 5
$ws.Range("E216").Value = "LP1912"
$ws.Range("A217").Value = "11:51:05"
$ws.Range("B217").Value = "12:21"
$ws.Range("C217").Value = "14_ABASTO"
$ws.Range("D217").Value = 30
$ws.Range("E217").Value = "LP1912"
$ws.Range("A224").Value = "11:51:05"
$ws.Range("B224").Value = "12:37"
$ws.Range("C224").Value = "27_EL RETIRO"
$ws.Range("D224").Value = 46
$ws.Range("E224").Value = "LP1912"
$ws.Range("A225").Value = "11:51:05"
$ws.Range("B225").Value = "12:37"
$ws.Range("C225").Value = "23_HERNANDEZ"
$ws.Range("D225").Value = 46
$ws.Range("E225").Value = "LP1912"
$ws.Range("A249").Value = "11:51:05"
$ws.Range("B249").Value = "13:20"
$ws.Range("C249").Value = "10_OLMOS"
$ws.Range("D249").Value = 89
$ws.Range("E249").Value = "LP1912"
$ws.Range("A250").Value = "11:23:54"
$ws.Range("B250").Value = "13:20"
$ws.Range("C250").Value = "26_HERNANDEZ"
$ws.Range("D250").Value = 117
$ws.Range("E250").Value = "LP1912"
$ws.Range("A337").Value = "15:51:40"
$ws.Range("B337").Value = "16:18"
$ws.Range("C337").Value = "16_SANTA ANA"
$ws.Range("D337").Value = 27
$ws.Range("E337").Value = "LP1912"
$ws.Range("A338").Value = "14:40:41"
$ws.Range("B338").Value = "16:19"
$ws.Range("C338").Value = "215C_EL PATO"
$ws.Range("D338").Value = 99
$ws.Range("E338").Value = "LP1912"
$ws.Range("A339").Value = "14:53:55"
$ws.Range("B339").Value = "16:20"
$ws.Range("C339").Value = "26_HERNANDEZ"
$ws.Range("D339").Value = 87
$ws.Range("E339").Value = "LP1912"
$ws.Range("A340").Value = "14:40:41"
$ws.Range("B340").Value = "16:21"
$ws.Range("C340").Value = "26_HERNANDEZ"
$ws.Range("D340").Value = 101
$ws.Range("E340").Value = "LP1912"
$ws.Range("A341").Value = "15:51:40"
$ws.Range("B341").Value = "16:29"
$ws.Range("C341").Value = "10_OLMOS"
$ws.Range("D341").Value = 38
$ws.Range("E341").Value = "LP1912"
$ws.Range("A342").Value = "14:40:41"
$ws.Range("B342").Value = "16:30"
$ws.Range("C342").Value = "15_ABASTO"
$ws.Range("D342").Value = 110
$ws.Range("E342").Value = "LP1912"
$ws.Range("A343").Value = "15:51:40"
$ws.Range("B343").Value = "16:34"
$ws.Range("C343").Value = "23_HERNANDEZ"
$ws.Range("D343").Value = 43
$ws.Range("E343").Value = "LP1912"
$ws.Range("A344").Value = "15:19:52"
$ws.Range("B344").Value = "16:36"
$ws.Range("C344").Value = "11_ETCHEVERRY"
$ws.Range("D344").Value = 77
$ws.Range("E344").Value = "LP1912"
$ws.Range("A345").Value = "15:19:52"
$ws.Range("B345").Value = "16:39"
$ws.Range("C345").Value = "17_ROMERO"
$ws.Range("D345").Value = 80
$ws.Range("E345").Value = "LP1912"
$ws.Range("A346").Value = "14:53:55"
$ws.Range("B346").Value = "16:42"
$ws.Range("C346").Value = "16_P MOR-SANTA ANA"
$ws.Range("D346").Value = 109
$ws.Range("E346").Value = "LP1912"
$ws.Range("A347").Value = "14:53:55"
$ws.Range("B347").Value = "16:42"
$ws.Range("C347").Value = "225_GOMEZ"
$ws.Range("D347").Value = 109
$ws.Range("E347").Value = "LP1912"
$ws.Range("A348").Value = "15:51:40"
$ws.Range("B348").Value = "16:43"
$ws.Range("C348").Value = "225_GOMEZ"
$ws.Range("D348").Value = 52
$ws.Range("E348").Value = "LP1912"
$ws.Range("A349").Value = "14:53:55"
$ws.Range("B349").Value = "16:48"
$ws.Range("C349").Value = "15_ABASTO"
$ws.Range("D349").Value = 115
$ws.Range("E349").Value = "LP1912"
$ws.Range("A350").Value = "15:51:40"
$ws.Range("B350").Value = "16:50"
$ws.Range("C350").Value = "14_ABASTO"
$ws.Range("D350").Value = 59
$ws.Range("E350").Value = "LP1912"
$ws.Range("A351").Value = "15:19:52"
$ws.Range("B351").Value = "16:56"
$ws.Range("C351").Value = "17_179 Y 38"
$ws.Range("D351").Value = 97
$ws.Range("E351").Value = "LP1912"
$ws.Range("A352").Value = "15:19:52"
$ws.Range("B352").Value = "17:04"
$ws.Range("C352").Value = "215A_EL PATO"
$ws.Range("D352").Value = 105
$ws.Range("E352").Value = "LP1912"
$ws.Range("A353").Value = "15:51:40"
$ws.Range("B353").Value = "17:04"
$ws.Range("C353").Value = "11_ETCHEVERRY"
$ws.Range("D353").Value = 73
$ws.Range("E353").Value = "LP1912"
$ws.Range("A354").Value = "15:51:40"
$ws.Range("B354").Value = "17:21"
$ws.Range("C354").Value = "26_HERNANDEZ"
$ws.Range("D354").Value = 90
$ws.Range("E354").Value = "LP1912"
$ws.Range("A355").Value = "15:51:40"
$ws.Range("B355").Value = "17:24"
$ws.Range("C355").Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Range("D355").Value = 93
$ws.Range("E355").Value = "LP1912"
$ws.Range("A356").Value = "15:51:40"
$ws.Range("B356").Value = "17:28"
$ws.Range("C356").Value = "14_ABASTO"
$ws.Range("D356").Value = 97
$ws.Range("E356").Value = "LP1912"
$ws.Range("A357").Value = "15:51:40"
$ws.Range("B357").Value = "17:36"
$ws.Range("C357").Value = "27_EL RETIRO"
$ws.Range("D357").Value = 105
$ws.Range("E357").Value = "LP1912"
$ws.Range("A358").Value = "15:51:40"
$ws.Range("B358").Value = "17:38"
$ws.Range("C358").Value = "17_ROMERO"
$ws.Range("D358").Value = 107
$ws.Range("E358").Value = "LP1912"
$ws.Range("A359").Value = "15:51:40"
$ws.Range("B359").Value = "17:40"
$ws.Range("C359").Value = "215B_EL PATO"
$ws.Range("D359").Value = 109
$ws.Range("E359").Value = "LP1912"
$ws.Range("A360").Value = "15:51:40"
$ws.Range("B360").Value = "17:50"
$ws.Range("C360").Value = "16_P MOR-167 Y 521"
$ws.Range("D360").Value = 119
$ws.Range("E360").Value = "LP1912"

$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 15:51:40"
$ws.Range("A3").Value = "Total filas: 38"
$ws.Range("A43").Value = "15:51:40"
$ws.Range("B43").Value = "17:40"
$ws.Range("C43").Value = "215B_EL PATO"
$ws.Range("D43").Value = 109
$ws.Range("E43").Value = "LP1912"

$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 15:51:40"
$ws.Range("A3").Value = "Total filas: 47"
$ws.Range("A52").Value = "15:51:40"
$ws.Range("B52").Value = "17:15"
$ws.Range("C52").Value = "215A_LA PLATA"
$ws.Range("D52").Value = 84
$ws.Range("E52").Value = "L6173"

